# Rename the original sheet and add the new CASH_FLOWS sheet, matching
# the target workbook layout: INCOME (sheetId 1), CASH_FLOWS (sheetId 2).
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "INCOME"

$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "CASH_FLOWS"

# ---- Header row (bold, bordered, centered - same look as INCOME!A1:D1) ----
$ws1.Range("A1:D1").Copy()
$ws2.Range("A1:E1").PasteSpecial(-4122)   # xlPasteFormats
$ws2.Application.CutCopyMode = $false

$ws2.Range("A1").Value = "Description"
$ws2.Range("B1").Value = "Value_1"
$ws2.Range("C1").Value = "Value_2"
$ws2.Range("D1").Value = "Value_3"
$ws2.Range("E1").Value = "Value_4"

# ---- Data rows. Values are stored as text (matching how the PDF-extracted
# INCOME sheet stores its numbers), so force text entry then drop back to
# the default "Normal" style to avoid leaving a stray quote-prefix format. ----
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws2.Range("A2") "Class A share repurchases"
Set-TextValue $ws2.Range("B2") "78"
Set-TextValue $ws2.Range("C2") "9316"
Set-TextValue $ws2.Range("D2") "73"
Set-TextValue $ws2.Range("E2") "11855"

Set-TextValue $ws2.Range("A3") "Class C share repurchases"
Set-TextValue $ws2.Range("B3") "450"
Set-TextValue $ws2.Range("C3") "52868"
Set-TextValue $ws2.Range("D3") "306"
Set-TextValue $ws2.Range("E3") "50192"

Set-TextValue $ws2.Range("A4") "Total share repurchases(1)"
Set-TextValue $ws2.Range("B4") "528"
Set-TextValue $ws2.Range("C4") "62184"
Set-TextValue $ws2.Range("D4") "379"
Set-TextValue $ws2.Range("E4") "62047"

# ---- Comments citing the PDF source page, on A1 of each sheet. ----
$ws1.Range("A1").AddComment("From page 54 of goog-10-k-2024.pdf")
$ws2.Range("A1").AddComment("From page 44 of goog-10-k-2024.pdf")

$ws1.Range("A1").Select() | Out-Null
